# Append a new scrape run (2026-01-05 06:44:47 JST):
#   - every existing case's "取得日時" (fetched-at) timestamp is refreshed
#     to the new run time
#   - a newly-found case ("クラウドウェア内製化推進...") is inserted as the
#     new row 7, pushing the former rows 7-8 down to rows 8-9
#   - the F-column hyperlinks are rebuilt so each URL cell keeps a working
#     link to its own (possibly shifted) target
#
# NB: `Range.Value` is a parameterized COM property in this runtime and a
# bare read of it doesn't resolve to the underlying scalar - use `Value2`
# whenever a cell's current value needs to be read back.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-05 06:44:47"

# --- shift the last two existing rows down by one to make room -------------
# (grab the old row 7/8 values first since row 7 is about to be overwritten)
$oldRow7_B = $ws.Cells.Item(7, 2).Value2
$oldRow7_C = $ws.Cells.Item(7, 3).Value2
$oldRow7_D = $ws.Cells.Item(7, 4).Value2
$oldRow7_E = $ws.Cells.Item(7, 5).Value2
$oldRow7_F = $ws.Cells.Item(7, 6).Value2
$oldRow7_G = $ws.Cells.Item(7, 7).Value2

$oldRow8_B = $ws.Cells.Item(8, 2).Value2
$oldRow8_C = $ws.Cells.Item(8, 3).Value2
$oldRow8_D = $ws.Cells.Item(8, 4).Value2
$oldRow8_E = $ws.Cells.Item(8, 5).Value2
$oldRow8_F = $ws.Cells.Item(8, 6).Value2
$oldRow8_G = $ws.Cells.Item(8, 7).Value2

# old row 8 -> new row 9
$ws.Cells.Item(9, 2).Value = $oldRow8_B
$ws.Cells.Item(9, 3).Value = $oldRow8_C
$ws.Cells.Item(9, 4).Value = $oldRow8_D
$ws.Cells.Item(9, 5).Value = $oldRow8_E
$ws.Cells.Item(9, 6).Value = $oldRow8_F
$ws.Cells.Item(9, 7).Value = $oldRow8_G

# old row 7 -> new row 8
$ws.Cells.Item(8, 2).Value = $oldRow7_B
$ws.Cells.Item(8, 3).Value = $oldRow7_C
$ws.Cells.Item(8, 4).Value = $oldRow7_D
$ws.Cells.Item(8, 5).Value = $oldRow7_E
$ws.Cells.Item(8, 6).Value = $oldRow7_F
$ws.Cells.Item(8, 7).Value = $oldRow7_G

# --- write the newly-scraped case into row 7 --------------------------------
$ws.Cells.Item(7, 2).Value = "【急募】クラウドウェア内製化推進のための技術サポート依頼"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5465210"
$ws.Cells.Item(7, 7).Value = 25

# --- refresh the "取得日時" column for every data row (2..9) ----------------
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- rebuild the F-column hyperlinks so they track the shifted rows --------
# (wiping via a cell-scoped Hyperlinks.Delete() clears the sheet's whole
# hyperlink collection in this runtime, so do it once up front, then re-add
# every link fresh in row order.)
$ws.Cells.Item(2, 6).Hyperlinks.Delete()

for ($r = 2; $r -le 9; $r++) {
    $url = $ws.Cells.Item($r, 6).Value2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url) | Out-Null
}

$ws.Range("A1").Select() | Out-Null
